$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from B1 into C1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header value for new column
$ws.Range("C1").Value = 2

# Data values for new column C (rows 2-6)
$ws.Range("C2").Value = -4.697124697347145
$ws.Range("C3").Value = -1.172924015787017
$ws.Range("C4").Value = -0.07027665786814449
$ws.Range("C5").Value = -0.4595607842740025
$ws.Range("C6").Value = -0.1300780636132118
